# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New columns AD, AE, AF are appended after the existing data (A:AC),
# with a header row and a constant record (86-76-0) repeated for every
# player row (2-52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border) from the last existing
# header cell (AC1) onto the three new header cells so they match the
# rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row.
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 86
    $ws.Range("AE$r").Value = 76
    $ws.Range("AF$r").Value = 0
}
